$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3747116666666667
$ws.Range("H2").Value = 1.124135
$ws.Range("I2").Value = 0.3914669751594584
$ws.Range("J2").Value = 0.3914669751594584
$ws.Range("M2").Value = 0.05661333333333334
$ws.Range("N2").Value = 0.16984
$ws.Range("O2").Value = 0.0204119846136133
$ws.Range("P2").Value = 0.02041198461361329
$ws.Range("Q2").Value = 0.0212136764888889
$ws.Range("R2").Value = 0.1909230884000001
$ws.Range("S2").Value = 0.007990617873692603
$ws.Range("T2").Value = 0.007990617873692601
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3747116666666667
$ws.Range("H3").Value = 1.124135
$ws.Range("I3").Value = 0.3914669751594584
$ws.Range("J3").Value = 0.3914669751594584
$ws.Range("N3").Value = 0.8341160000000001
$ws.Range("O3").Value = 0.100247073468963
$ws.Range("P3").Value = 0.1002470734689629
$ws.Range("Q3").Value = 0.1041843321844445
$ws.Range("R3").Value = 0.9376589896600002
$ws.Range("S3").Value = 0.03924341861948293
$ws.Range("T3").Value = 0.03924341861948292
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3747116666666667
$ws.Range("H4").Value = 1.124135
$ws.Range("I4").Value = 0.3914669751594584
$ws.Range("J4").Value = 0.3914669751594584
$ws.Range("M4").Value = 2.438882
$ws.Range("N4").Value = 7.316646
$ws.Range("O4").Value = 0.8793409419174237
$ws.Range("P4").Value = 0.8793409419174237
$ws.Range("Q4").Value = 0.9138775390233335
$ws.Range("R4").Value = 8.224897851210001
$ws.Range("S4").Value = 0.3442329386662829
$ws.Range("T4").Value = 0.3442329386662829
$ws.Range("I5").Value = 0.4195022558883632
$ws.Range("J5").Value = 0.4195022558883631
$ws.Range("M5").Value = 0.05661333333333334
$ws.Range("N5").Value = 0.16984
$ws.Range("O5").Value = 0.0204119846136133
$ws.Range("P5").Value = 0.02041198461361329
$ws.Range("Q5").Value = 0.02273291416000001
$ws.Range("R5").Value = 0.20459622744
$ws.Range("S5").Value = 0.008562873592569336
$ws.Range("T5").Value = 0.008562873592569334
$ws.Range("I6").Value = 0.4195022558883632
$ws.Range("J6").Value = 0.4195022558883631
$ws.Range("N6").Value = 0.8341160000000001
$ws.Range("O6").Value = 0.100247073468963
$ws.Range("P6").Value = 0.1002470734689629
$ws.Range("S6").Value = 0.04205387346643644
$ws.Range("T6").Value = 0.04205387346643643
$ws.Range("I7").Value = 0.4195022558883632
$ws.Range("J7").Value = 0.4195022558883631
$ws.Range("M7").Value = 2.438882
$ws.Range("N7").Value = 7.316646
$ws.Range("O7").Value = 0.8793409419174237
$ws.Range("P7").Value = 0.8793409419174237
$ws.Range("Q7").Value = 0.9793257504540001
$ws.Range("R7").Value = 8.813931754086001
$ws.Range("S7").Value = 0.3688855088293574
$ws.Range("T7").Value = 0.3688855088293573
$ws.Range("G8").Value = 0.18094
$ws.Range("H8").Value = 0.54282
$ws.Range("I8").Value = 0.1890307689521785
$ws.Range("J8").Value = 0.1890307689521785
$ws.Range("M8").Value = 0.05661333333333334
$ws.Range("N8").Value = 0.16984
$ws.Range("O8").Value = 0.0204119846136133
$ws.Range("P8").Value = 0.02041198461361329
$ws.Range("Q8").Value = 0.01024361653333333
$ws.Range("R8").Value = 0.09219254880000001
$ws.Range("S8").Value = 0.003858493147351357
$ws.Range("T8").Value = 0.003858493147351357
$ws.Range("G9").Value = 0.18094
$ws.Range("H9").Value = 0.54282
$ws.Range("I9").Value = 0.1890307689521785
$ws.Range("J9").Value = 0.1890307689521785
$ws.Range("N9").Value = 0.8341160000000001
$ws.Range("O9").Value = 0.100247073468963
$ws.Range("P9").Value = 0.1002470734689629
$ws.Range("Q9").Value = 0.05030831634666667
$ws.Range("R9").Value = 0.45277484712
$ws.Range("S9").Value = 0.0189497813830436
$ws.Range("T9").Value = 0.0189497813830436
$ws.Range("G10").Value = 0.18094
$ws.Range("H10").Value = 0.54282
$ws.Range("I10").Value = 0.1890307689521785
$ws.Range("J10").Value = 0.1890307689521785
$ws.Range("M10").Value = 2.438882
$ws.Range("N10").Value = 7.316646
$ws.Range("O10").Value = 0.8793409419174237
$ws.Range("P10").Value = 0.8793409419174237
$ws.Range("Q10").Value = 0.44129130908
$ws.Range("R10").Value = 3.97162178172
$ws.Range("S10").Value = 0.1662224944217835
$ws.Range("T10").Value = 0.1662224944217835
